# [IMP] Data with date range 2024
# Adds a new block of date_range rows (fiscal year, 12 months, 4 quarters)
# for year 2024, following the same pattern as the existing 2020-2023 blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the formatting (styles) of the last existing data row (69)
#        down onto the new rows (70-86) so the new cells use the same
#        styles (s=2 for text/number cells, s=3 for date cells) as the
#        rest of the table instead of Excel's default style. ---
$ws.Range("A69:H69").Copy()
$ws.Range("A70:H86").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Fill in the new rows with data -------------------------------

# Row 70: z0bug.2024 (fiscal year)
$ws.Range("A70").Value = "z0bug.2024"
$ws.Range("B70").Value = "z0bug.fiscal"
$ws.Range("C70").Value = 2024
$ws.Range("D70").Value = 45292
$ws.Range("E70").Value = 45657
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 1

# Row 71: z0bug.2024-01 (January)
$ws.Range("A71").Value = "z0bug.2024-01"
$ws.Range("B71").Value = "z0bug.monthly"
$ws.Range("C71").Value = "2024/01"
$ws.Range("D71").Value = 45292
$ws.Range("E71").Value = 45322
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 1

# Row 72: z0bug.2024-02 (February)
$ws.Range("A72").Value = "z0bug.2024-02"
$ws.Range("B72").Value = "z0bug.monthly"
$ws.Range("C72").Value = "2024/02"
$ws.Range("D72").Value = 45323
$ws.Range("E72").Value = 45350
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 1

# Row 73: z0bug.2024-03 (March)
$ws.Range("A73").Value = "z0bug.2024-03"
$ws.Range("B73").Value = "z0bug.monthly"
$ws.Range("C73").Value = "2024/03"
$ws.Range("D73").Value = 45352
$ws.Range("E73").Value = 45382
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 1

# Row 74: z0bug.2024-04 (April)
$ws.Range("A74").Value = "z0bug.2024-04"
$ws.Range("B74").Value = "z0bug.monthly"
$ws.Range("C74").Value = "2024/04"
$ws.Range("D74").Value = 45383
$ws.Range("E74").Value = 45412
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 1

# Row 75: z0bug.2024-05 (May)
$ws.Range("A75").Value = "z0bug.2024-05"
$ws.Range("B75").Value = "z0bug.monthly"
$ws.Range("C75").Value = "2024/05"
$ws.Range("D75").Value = 45413
$ws.Range("E75").Value = 45443
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 1

# Row 76: z0bug.2024-06 (June)
$ws.Range("A76").Value = "z0bug.2024-06"
$ws.Range("B76").Value = "z0bug.monthly"
$ws.Range("C76").Value = "2024/06"
$ws.Range("D76").Value = 45444
$ws.Range("E76").Value = 45473
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 1

# Row 77: z0bug.2024-07 (July)
$ws.Range("A77").Value = "z0bug.2024-07"
$ws.Range("B77").Value = "z0bug.monthly"
$ws.Range("C77").Value = "2024/07"
$ws.Range("D77").Value = 45474
$ws.Range("E77").Value = 45504
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 1

# Row 78: z0bug.2024-08 (August)
$ws.Range("A78").Value = "z0bug.2024-08"
$ws.Range("B78").Value = "z0bug.monthly"
$ws.Range("C78").Value = "2024/08"
$ws.Range("D78").Value = 45505
$ws.Range("E78").Value = 45535
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 1

# Row 79: z0bug.2024-09 (September)
$ws.Range("A79").Value = "z0bug.2024-09"
$ws.Range("B79").Value = "z0bug.monthly"
$ws.Range("C79").Value = "2024/09"
$ws.Range("D79").Value = 45536
$ws.Range("E79").Value = 45565
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 1

# Row 80: z0bug.2024-10 (October)
$ws.Range("A80").Value = "z0bug.2024-10"
$ws.Range("B80").Value = "z0bug.monthly"
$ws.Range("C80").Value = "2024/10"
$ws.Range("D80").Value = 45566
$ws.Range("E80").Value = 45596
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 1

# Row 81: z0bug.2024-11 (November)
$ws.Range("A81").Value = "z0bug.2024-11"
$ws.Range("B81").Value = "z0bug.monthly"
$ws.Range("C81").Value = "2024/11"
$ws.Range("D81").Value = 45597
$ws.Range("E81").Value = 45626
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 1

# Row 82: z0bug.2024-12 (December)
$ws.Range("A82").Value = "z0bug.2024-12"
$ws.Range("B82").Value = "z0bug.monthly"
$ws.Range("C82").Value = "2024/12"
$ws.Range("D82").Value = 45627
$ws.Range("E82").Value = 45657
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 1

# Row 83: z0bug.2024-Q1
$ws.Range("A83").Value = "z0bug.2024-Q1"
$ws.Range("B83").Value = "z0bug.quarter"
$ws.Range("C83").Value = "2024/T1"
$ws.Range("D83").Value = 45292
$ws.Range("E83").Value = 45382
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 1

# Row 84: z0bug.2024-Q2
$ws.Range("A84").Value = "z0bug.2024-Q2"
$ws.Range("B84").Value = "z0bug.quarter"
$ws.Range("C84").Value = "2024/T2"
$ws.Range("D84").Value = 45383
$ws.Range("E84").Value = 45473
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 1

# Row 85: z0bug.2024-Q3
$ws.Range("A85").Value = "z0bug.2024-Q3"
$ws.Range("B85").Value = "z0bug.quarter"
$ws.Range("C85").Value = "2024/T3"
$ws.Range("D85").Value = 45474
$ws.Range("E85").Value = 45565
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 1

# Row 86: z0bug.2024-Q4
$ws.Range("A86").Value = "z0bug.2024-Q4"
$ws.Range("B86").Value = "z0bug.quarter"
$ws.Range("C86").Value = "2024/T4"
$ws.Range("D86").Value = 45566
$ws.Range("E86").Value = 45657
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 1

# --- 3. Match row heights of the new rows to the rest of the table ---
for ($r = 70; $r -le 86; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# --- 4. Update the view: scroll down and select the newly added data -
$ws.Range("A61").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1

$selRange = $ws.Range("A70:A86,C70:E86")
$selRange.Select()
